# Aggiornamento dati COVID comune di Formigine al 23 agosto 2021
# Appends daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
# for rows 344-357 (dates 2021-08-10 .. 2021-08-23) to the existing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44418, 6, 39, 113.4842576965606),
    @(44419, 0, 38, 110.5744049351103),
    @(44420, 6, 40, 116.3941104580108),
    @(44421, 7, 42, 122.2138159809114),
    @(44422, 8, 41, 119.3039632194611),
    @(44423, 7, 38, 110.5744049351103),
    @(44424, 5, 39, 113.4842576965606),
    @(44425, 7, 40, 116.3941104580108),
    @(44426, 1, 41, 119.3039632194611),
    @(44427, 5, 40, 116.3941104580108),
    @(44428, 7, 40, 116.3941104580108),
    @(44429, 0, 32, 93.11528836640865),
    @(44430, 5, 30, 87.29558284350811),
    @(44431, 2, 27, 78.56602455915731)
)

$startRow = 344
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Copy column A's date formatting (style index) from the last pre-existing row
# down across the newly appended rows, matching the original sheet's styling.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A$startRow`:A357").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
